# Update UNLV_A team-specific matrix with recomputed percentages
# (games pulled March 7). Cell addresses/values below correspond 1:1
# to the rows/columns changed in the source XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value = 0.2135231316725979    # B2
$ws.Cells.Item(2, 3).Value = 0.5338078291814946    # C2
$ws.Cells.Item(2, 10).Value = 0.01067615658362989  # J2
$ws.Cells.Item(2, 16).Value = 0.1423487544483986   # P2
$ws.Cells.Item(2, 19).Value = 0.099644128113879    # S2

# Row 3
$ws.Cells.Item(3, 3).Value = 0.01935483870967742   # C3
$ws.Cells.Item(3, 10).Value = 0.01290322580645161  # J3
$ws.Cells.Item(3, 16).Value = 0.7870967741935484   # P3
$ws.Cells.Item(3, 19).Value = 0.1806451612903226   # S3

# Row 4
$ws.Cells.Item(4, 10).Value = 0.01724137931034483  # J4
$ws.Cells.Item(4, 16).Value = 0.7586206896551724   # P4
$ws.Cells.Item(4, 19).Value = 0.2241379310344828   # S4

# Row 5
$ws.Cells.Item(5, 15).Value = 0.3333333333333333   # O5
$ws.Cells.Item(5, 16).Value = 0.3333333333333333   # P5
$ws.Cells.Item(5, 19).Value = 0.3333333333333333   # S5

# Row 6
$ws.Cells.Item(6, 2).Value = 0.04524886877828054   # B6
$ws.Cells.Item(6, 4).Value = 0.004524886877828055  # D6
$ws.Cells.Item(6, 5).Value = 0.004524886877828055  # E6
$ws.Cells.Item(6, 6).Value = 0.04072398190045249   # F6
$ws.Cells.Item(6, 10).Value = 0.2895927601809955   # J6
$ws.Cells.Item(6, 15).Value = 0.004524886877828055 # O6
$ws.Cells.Item(6, 17).Value = 0.16289592760181     # Q6
$ws.Cells.Item(6, 18).Value = 0.04977375565610859  # R6
$ws.Cells.Item(6, 19).Value = 0.3981900452488688   # S6

# Row 7
$ws.Cells.Item(7, 2).Value = 0.0989010989010989    # B7
$ws.Cells.Item(7, 4).Value = 0.04945054945054945   # D7
$ws.Cells.Item(7, 6).Value = 0.05494505494505494   # F7
$ws.Cells.Item(7, 10).Value = 0.1208791208791209   # J7
$ws.Cells.Item(7, 15).Value = 0.03296703296703297  # O7
$ws.Cells.Item(7, 17).Value = 0.2032967032967033   # Q7
$ws.Cells.Item(7, 18).Value = 0.07142857142857142  # R7
$ws.Cells.Item(7, 19).Value = 0.3681318681318682   # S7

# Row 8
$ws.Cells.Item(8, 2).Value = 0.08888888888888889   # B8
$ws.Cells.Item(8, 4).Value = 0.01728395061728395   # D8
$ws.Cells.Item(8, 6).Value = 0.05925925925925926   # F8
$ws.Cells.Item(8, 10).Value = 0.1234567901234568   # J8
$ws.Cells.Item(8, 15).Value = 0.01234567901234568  # O8
$ws.Cells.Item(8, 17).Value = 0.1703703703703704   # Q8
$ws.Cells.Item(8, 18).Value = 0.1012345679012346   # R8
$ws.Cells.Item(8, 19).Value = 0.4271604938271605   # S8

# Row 9
$ws.Cells.Item(9, 2).Value = 0.1067415730337079    # B9
$ws.Cells.Item(9, 4).Value = 0.005617977528089887  # D9
$ws.Cells.Item(9, 5).Value = 0.005617977528089887  # E9
$ws.Cells.Item(9, 6).Value = 0.07865168539325842   # F9
$ws.Cells.Item(9, 10).Value = 0.09550561797752809  # J9
$ws.Cells.Item(9, 15).Value = 0.01123595505617977  # O9
$ws.Cells.Item(9, 17).Value = 0.2191011235955056   # Q9
$ws.Cells.Item(9, 18).Value = 0.07865168539325842  # R9
$ws.Cells.Item(9, 19).Value = 0.398876404494382    # S9

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1129431162407255   # B10
$ws.Cells.Item(10, 4).Value = 0.0313272877164056   # D10
$ws.Cells.Item(10, 5).Value = 0.0008244023083264633 # E10
$ws.Cells.Item(10, 6).Value = 0.07007419620774938  # F10
$ws.Cells.Item(10, 10).Value = 0.1294311624072547  # J10
$ws.Cells.Item(10, 15).Value = 0.01483924154987634 # O10
$ws.Cells.Item(10, 17).Value = 0.2308326463314097  # Q10
$ws.Cells.Item(10, 18).Value = 0.07749381698268755 # R10
$ws.Cells.Item(10, 19).Value = 0.3322341302555647  # S10

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1085271317829457   # G11
$ws.Cells.Item(11, 10).Value = 0.09302325581395349 # J11
$ws.Cells.Item(11, 11).Value = 0.1627906976744186  # K11
$ws.Cells.Item(11, 12).Value = 0.6085271317829457  # L11
$ws.Cells.Item(11, 19).Value = 0.02713178294573643 # S11

# Row 12
$ws.Cells.Item(12, 7).Value = 0.7939393939393939   # G12
$ws.Cells.Item(12, 10).Value = 0.1333333333333333  # J12
$ws.Cells.Item(12, 11).Value = 0.01818181818181818 # K12
$ws.Cells.Item(12, 12).Value = 0.0303030303030303  # L12
$ws.Cells.Item(12, 19).Value = 0.02424242424242424 # S12

# Row 15
$ws.Cells.Item(15, 6).Value = 0.01834862385321101  # F15
$ws.Cells.Item(15, 8).Value = 0.1422018348623853   # H15
$ws.Cells.Item(15, 9).Value = 0.05504587155963303  # I15
$ws.Cells.Item(15, 10).Value = 0.4082568807339449  # J15
$ws.Cells.Item(15, 11).Value = 0.04587155963302753 # K15
$ws.Cells.Item(15, 13).Value = 0.009174311926605505 # M15
$ws.Cells.Item(15, 15).Value = 0.05045871559633028 # O15
$ws.Cells.Item(15, 19).Value = 0.2706422018348624  # S15

# Row 16
$ws.Cells.Item(16, 6).Value = 0.02985074626865672  # F16
$ws.Cells.Item(16, 8).Value = 0.1293532338308458   # H16
$ws.Cells.Item(16, 9).Value = 0.09950248756218906  # I16
$ws.Cells.Item(16, 10).Value = 0.3880597014925373  # J16
$ws.Cells.Item(16, 11).Value = 0.09950248756218906 # K16
$ws.Cells.Item(16, 13).Value = 0.01990049751243781 # M16
$ws.Cells.Item(16, 14).Value = 0.004975124378109453 # N16
$ws.Cells.Item(16, 15).Value = 0.05970149253731343 # O16
$ws.Cells.Item(16, 19).Value = 0.1691542288557214  # S16

# Row 17
$ws.Cells.Item(17, 6).Value = 0.01746724890829694  # F17
$ws.Cells.Item(17, 8).Value = 0.1877729257641921   # H17
$ws.Cells.Item(17, 9).Value = 0.08078602620087336  # I17
$ws.Cells.Item(17, 10).Value = 0.4213973799126637  # J17
$ws.Cells.Item(17, 11).Value = 0.08078602620087336 # K17
$ws.Cells.Item(17, 13).Value = 0.01965065502183406 # M17
$ws.Cells.Item(17, 15).Value = 0.07423580786026202 # O17
$ws.Cells.Item(17, 19).Value = 0.1179039301310044  # S17

# Row 18
$ws.Cells.Item(18, 6).Value = 0.03508771929824561  # F18
$ws.Cells.Item(18, 8).Value = 0.1988304093567251   # H18
$ws.Cells.Item(18, 9).Value = 0.0935672514619883   # I18
$ws.Cells.Item(18, 10).Value = 0.4035087719298245  # J18
$ws.Cells.Item(18, 11).Value = 0.08187134502923976 # K18
$ws.Cells.Item(18, 13).Value = 0.01754385964912281 # M18
$ws.Cells.Item(18, 15).Value = 0.05847953216374269 # O18
$ws.Cells.Item(18, 19).Value = 0.1111111111111111  # S18

# Row 19
$ws.Cells.Item(19, 6).Value = 0.0188034188034188   # F19
$ws.Cells.Item(19, 8).Value = 0.1974358974358974   # H19
$ws.Cells.Item(19, 9).Value = 0.07777777777777778  # I19
$ws.Cells.Item(19, 10).Value = 0.3649572649572649  # J19
$ws.Cells.Item(19, 11).Value = 0.1111111111111111  # K19
$ws.Cells.Item(19, 13).Value = 0.02136752136752137 # M19
$ws.Cells.Item(19, 15).Value = 0.08632478632478632 # O19
$ws.Cells.Item(19, 19).Value = 0.1222222222222222  # S19
